$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(1).Copy()
$ws.Rows.Item(3).PasteSpecial(-4163)
$ws.Rows.Item(3).ClearFormats()
$v = $ws.Cells.Item(3,1).Value()
$t = $ws.Cells.Item(3,1).Value().GetType().Name
Write-Host "A3 = [$v] type=$t"
